# Adds a new "Code of Conduct" section to the tender document, right
# between the "Pricing" section ("$200,000") and the "Quality" section:
#   - fills the existing empty paragraph that follows the pricing
#     paragraph with the heading text "Code of Conduct"
#   - inserts a brand new paragraph after it containing the body text
#     "The company has a strict code of conduct policy which aligns
#      with that of the Victorian Public Sector."

$d = $word.ActiveDocument

# Find the "$200,000" pricing line, then step to the (currently empty)
# paragraph right after it - that's where the new section goes.
$findRange = $d.Content
$found = $findRange.Find.Execute("`$200,000", $true, $false, $false, $false, `
                                  $false, $true, 1, $false, "", 0)
if (-not $found) {
    throw "Could not locate the '$200,000' pricing paragraph"
}

$pricingPara = $findRange.Paragraphs(1)
$target = $pricingPara.Next()

# Fill in the heading text for the paragraph, matching the surrounding
# runs' en-AU language formatting.
$target.Range.Text = "Code of Conduct"
$target.Range.LanguageID = "en-AU"

# Insert a fresh paragraph right after it and fill it with the body text.
$target.Range.InsertParagraphAfter()

$body = $target.Next()
$body.Range.Text = "The company has a strict code of conduct policy which aligns with that of the Victorian Public Sector."
$body.Range.LanguageID = "en-AU"
